# Replace each arithmetic expression in the practice-sheet table with a new
# expression, one Find/Replace per cell value. Expressions are looked up and
# replaced using Range.Find over the whole document content; each call only
# replaces a single occurrence (wdReplaceOne) to keep a 1:1 mapping between
# old and new values even though the runtime's Find ignores Range scoping.
#
# NOTE: "22+17=" is processed before "2+17=" because "2+17=" is a literal
# substring of "22+17="; replacing the shorter pattern first would corrupt
# the "22+17=" cell (it would match inside it) before that cell gets its own
# correct replacement.
$d = $word.ActiveDocument

$replacements = @(
    @('38-21=', '62+31='),
    @('6+19=', '94-0='),
    @('46+27=', '15+29='),
    @('31-16=', '55-52='),
    @('80-14=', '42+19='),
    @('66-0=', '87-14='),
    @('88-41=', '7+33='),
    @('35+46=', '68-7='),
    @('31+39=', '45-25='),
    @('75-49=', '23+13='),
    @('3+90=', '6+79='),
    @('46-22=', '24-8='),
    @('78-24=', '17+73='),
    @('54-37=', '2+84='),
    @('20+3=', '32-24='),
    @('1+72=', '19+35='),
    @('93-92=', '17-12='),
    @('23-9=', '34+7='),
    @('24+62=', '97-22='),
    @('30+18=', '72-23='),
    @('53+2=', '63-52='),
    @('5+74=', '85-41='),
    @('1+76=', '73-49='),
    @('80+8=', '44+44='),
    @('85-63=', '97-36='),
    @('72-22=', '46+30='),
    @('15+47=', '68+14='),
    @('38+53=', '11+39='),
    @('22+17=', '65-1='),
    @('2+17=', '82-32='),
    @('10+28=', '86-69='),
    @('17+37=', '41+2='),
    @('62+28=', '77+12='),
    @('78+2=', '18+34='),
    @('87-42=', '38-25='),
    @('33-26=', '67-44='),
    @('87-46=', '77+19='),
    @('57+15=', '26+7='),
    @('69-58=', '83-54='),
    @('40-2=', '38+4='),
    @('6+45=', '24-11='),
    @('45+37=', '24+4='),
    @('25+8=', '63+14='),
    @('11-7=', '13+60='),
    @('14-6=', '93-56='),
    @('98-98=', '79-74='),
    @('64-64=', '4+34='),
    @('60+35=', '75-66='),
    @('62-41=', '25+64='),
    @('85+5=', '17+31='),
    @('63+3=', '1+7='),
    @('52-49=', '92+1='),
    @('65-57=', '43+50='),
    @('90+4=', '32+1='),
    @('44+23=', '50-15='),
    @('41+8=', '77-44='),
    @('92-82=', '92-5='),
    @('55-10=', '10+41='),
    @('99-74=', '26+3='),
    @('17+43=', '13+61='),
    @('37-32=', '17+80='),
    @('79+8=', '91-47='),
    @('87-21=', '15+39='),
    @('65-64=', '70+2='),
    @('55-51=', '60+23='),
    @('21+10=', '60-11='),
    @('36+56=', '5+9='),
    @('50-21=', '46-42='),
    @('17+59=', '38-29='),
    @('83-17=', '75+16='),
    @('46-38=', '39+15='),
    @('89-37=', '81+13='),
    @('65-0=', '23+68='),
    @('64-54=', '34+18='),
    @('65-29=', '99-62='),
    @('50-9=', '90-1='),
    @('10+31=', '84-84='),
    @('97-6=', '55-41='),
    @('2+36=', '86-83='),
    @('67-25=', '23+73='),
    @('88-20=', '21+77='),
    @('81-80=', '49-8='),
    @('27+46=', '95+3='),
    @('62-56=', '7+40='),
    @('92-34=', '62-29='),
    @('3+94=', '48-7='),
    @('1+14=', '79-71='),
    @('30+60=', '47-5='),
    @('38+56=', '18+77='),
    @('16+16=', '34+56='),
    @('96-56=', '16+63='),
    @('93-23=', '66+23='),
    @('71-30=', '49+12='),
    @('87-51=', '64+3='),
    @('97-88=', '79-16='),
    @('5+48=', '33+36='),
    @('69-2=', '74-33='),
    @('58+10=', '94-30='),
    @('69-38=', '79-67='),
    @('71+9=', '12+16=')
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}
